# atualizacao matriz de risco
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $ok = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND: $old"
    }
}

# --- Indicadores / prazos (dias) updates ---
Replace-Text "Concluído (97 dias)" "Concluído (110 dias)"
Replace-Text "Assinatura Contrato (106 dias)" "Assinatura Contrato (119 dias)"
Replace-Text "Assinatura Contrato (14 dias)" "Assinatura Contrato (27 dias)"
Replace-Text "Concluído (22 dias)" "Concluído (35 dias)"
Replace-Text "Assinatura Contrato (49 dias)" "Assinatura Contrato (62 dias)"
Replace-Text "Total de dias 708" "Total de dias 773"

# --- Matriz de risco / numero do processo updates ---
Replace-Text "785810/2024-010/00" "785810/2024-055/00"
Replace-Text "785810/2024-011/00" "785810/2024-056/00"
Replace-Text "785810/2024-012/00" "785810/2024-057/00"
Replace-Text "785810/2024-013/00" "785810/2024-058/00"
Replace-Text "785810/2024-014/00" "785810/2024-059/00"
Replace-Text "785810/2024-015/00" "785810/2024-060/00"
Replace-Text "785810/2024-016/00" "785810/2024-061/00"
Replace-Text "785810/2024-017/00" "785810/2024-062/00"
Replace-Text "785810/2024-018/00" "785810/2024-063/00"
Replace-Text "785810/2024-019/00" "785810/2024-064/00"
Replace-Text "785810/2024-020/00" "785810/2024-065/00"
Replace-Text "785810/2024-021/00" "785810/2024-066/00"

Write-Output "Done."
